$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder ("datetimeFigureOut" field) on the slide master and on
#    every slide layout: 04.02.2026 -> 06.02.2026
# ---------------------------------------------------------------------------
$oldDate = "04.02.2026"
$newDate = "06.02.2026"

function Update-DateField($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout
}

# ---------------------------------------------------------------------------
# 2) Slide 5 ("TextBox 9"): merge the trailing " " run with the following
#    "192.168.11.1" run into a single " 192.168.11.1" run.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$shp5 = $slide5.Shapes.Item(2)
$tr5 = $shp5.TextFrame.TextRange
$full5 = $tr5.Text
$needle5 = " 192.168.11.1"
$pos5 = $full5.IndexOf($needle5)
if ($pos5 -ge 0) {
    $sub5 = $tr5.Characters($pos5 + 1, $needle5.Length)
    $sub5.Text = $needle5
}

# ---------------------------------------------------------------------------
# 3) Slide 7 ("TextBox 10"): "led" -> "led " (add trailing space before "на").
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$shp7 = $slide7.Shapes.Item(4)
$tr7 = $shp7.TextFrame.TextRange
$full7 = $tr7.Text
$needle7 = "led"
$pos7 = $full7.IndexOf($needle7)
if ($pos7 -ge 0) {
    $sub7 = $tr7.Characters($pos7 + 1, $needle7.Length)
    $sub7.Text = "led "
}
